$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 7721
    $ws.Range("F4").Value = 530
    $ws.Range("F9").Value = 5954
    $ws.Range("F10").Value = 153
    $ws.Range("F13").Value = 1815
    $ws.Range("F14").Value = 1329
    $ws.Range("F15").Value = 285
    $ws.Range("F16").Value = 619
    $ws.Range("F17").Value = 138
    $ws.Range("F19").Value = 66
}
